$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.390.65"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "3.504.60"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'587.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'136.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("D7").Value = "3.505.67"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").Value = "'7.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "'0.377"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.94%  "
$ws.Range("D13").Value = "4.105.21"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "3.507.11"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "64.357.96"
$ws.Range("D18").Value = "'25.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.79%  "
$ws.Range("D19").Value = "'9.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "'13.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").Value = "'5.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "'384.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D24").Value = "3.644.95"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").Value = "'74.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'5.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("D29").Value = "'1.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "'7.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "'8.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").Value = "3.525.34"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("D38").Value = "'5.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").Value = "'6.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("D41").Value = "'163.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.58%  "
$ws.Range("E42").Value = "  -2.45%  "
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").Value = "'26.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'41.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "'1.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "2.479.00"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("E51").Value = "  -1.52%  "
